# Models_Information.xlsx — "add a few files"
#
# Content change: every "hardware_info" (column C) cell that reads
# "tested on MP34DT05" gets a sample-rate qualifier appended, becoming
# "tested on MP34DT05 20000 sps". Rows 7, 8, 9, 10, 11, 12, 17 and 20
# (the Micro_Speech_* and Voice_Turn models) carry that value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(7, 8, 9, 10, 11, 12, 17, 20)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Text -eq "tested on MP34DT05") {
        $cell.Value = "tested on MP34DT05 20000 sps"
    }
}

# Leave the workbook/view the way the author left it when they saved:
# scrolled back to column A, with the last-touched cell (C20) selected.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("C20").Select()
